$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.987569333333333
$ws.Range("H2").Value = 5.962707999999999
$ws.Range("I2").Value = 0.1241595834663642
$ws.Range("J2").Value = 0.1241595834663642
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 26.83081766666667
$ws.Range("N2").Value = 80.49245300000001
$ws.Range("O2").Value = 0.5916656861001716
$ws.Range("P2").Value = 0.5916656861001716
$ws.Range("Q2").Value = 53.32811038252489
$ws.Range("R2").Value = 479.952993442724
$ws.Range("S2").Value = 0.0734609651375379
$ws.Range("T2").Value = 0.07346096513753791

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.987569333333333
$ws.Range("H3").Value = 5.962707999999999
$ws.Range("I3").Value = 0.1241595834663642
$ws.Range("J3").Value = 0.1241595834663642
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.21969166666667
$ws.Range("N3").Value = 30.659075
$ws.Range("O3").Value = 0.2253617819930474
$ws.Range("P3").Value = 0.2253617819930474
$ws.Range("Q3").Value = 20.31234575278889
$ws.Range("R3").Value = 182.8111117751
$ws.Range("S3").Value = 0.02798082498149435
$ws.Range("T3").Value = 0.02798082498149435

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.987569333333333
$ws.Range("H4").Value = 5.962707999999999
$ws.Range("I4").Value = 0.1241595834663642
$ws.Range("J4").Value = 0.1241595834663642
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.297426666666667
$ws.Range("N4").Value = 24.89228
$ws.Range("O4").Value = 0.1829725319067811
$ws.Range("P4").Value = 0.1829725319067811
$ws.Range("Q4").Value = 16.49171078824889
$ws.Range("R4").Value = 148.42539709424
$ws.Range("S4").Value = 0.02271779334733197
$ws.Range("T4").Value = 0.02271779334733198

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.481595333333334
$ws.Range("H5").Value = 22.444786
$ws.Range("I5").Value = 0.4673606825542495
$ws.Range("J5").Value = 0.4673606825542495
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 26.83081766666667
$ws.Range("N5").Value = 80.49245300000001
$ws.Range("O5").Value = 0.5916656861001716
$ws.Range("P5").Value = 0.5916656861001716
$ws.Range("Q5").Value = 200.7373202444509
$ws.Range("R5").Value = 1806.635882200058
$ws.Range("S5").Value = 0.2765212788997045
$ws.Range("T5").Value = 0.2765212788997045

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 7.481595333333334
$ws.Range("H6").Value = 22.444786
$ws.Range("I6").Value = 0.4673606825542495
$ws.Range("J6").Value = 0.4673606825542495
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.21969166666667
$ws.Range("N6").Value = 30.659075
$ws.Range("O6").Value = 0.2253617819930474
$ws.Range("P6").Value = 0.2253617819930474
$ws.Range("Q6").Value = 76.45959748143891
$ws.Range("R6").Value = 688.13637733295
$ws.Range("S6").Value = 0.1053252362539126
$ws.Range("T6").Value = 0.1053252362539126

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.481595333333334
$ws.Range("H7").Value = 22.444786
$ws.Range("I7").Value = 0.4673606825542495
$ws.Range("J7").Value = 0.4673606825542495
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.297426666666667
$ws.Range("N7").Value = 24.89228
$ws.Range("O7").Value = 0.1829725319067811
$ws.Range("P7").Value = 0.1829725319067811
$ws.Range("Q7").Value = 62.07798862800889
$ws.Range("R7").Value = 558.70189765208
$ws.Range("S7").Value = 0.08551416740063238
$ws.Range("T7").Value = 0.08551416740063239

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.539018333333334
$ws.Range("H8").Value = 19.617055
$ws.Range("I8").Value = 0.4084797339793862
$ws.Range("J8").Value = 0.4084797339793863
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 26.83081766666667
$ws.Range("N8").Value = 80.49245300000001
$ws.Range("O8").Value = 0.5916656861001716
$ws.Range("P8").Value = 0.5916656861001716
$ws.Range("Q8").Value = 175.4472086206573
$ws.Range("R8").Value = 1579.024877585915
$ws.Range("S8").Value = 0.2416834420629291
$ws.Range("T8").Value = 0.2416834420629291

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.539018333333334
$ws.Range("H9").Value = 19.617055
$ws.Range("I9").Value = 0.4084797339793862
$ws.Range("J9").Value = 0.4084797339793863
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 10.21969166666667
$ws.Range("N9").Value = 30.659075
$ws.Range("O9").Value = 0.2253617819930474
$ws.Range("P9").Value = 0.2253617819930474
$ws.Range("Q9").Value = 66.82675116934723
$ws.Range("R9").Value = 601.440760524125
$ws.Range("S9").Value = 0.09205572075764043
$ws.Range("T9").Value = 0.09205572075764044

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.539018333333334
$ws.Range("H10").Value = 19.617055
$ws.Range("I10").Value = 0.4084797339793862
$ws.Range("J10").Value = 0.4084797339793863
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.297426666666667
$ws.Range("N10").Value = 24.89228
$ws.Range("O10").Value = 0.1829725319067811
$ws.Range("P10").Value = 0.1829725319067811
$ws.Range("Q10").Value = 54.25702509282223
$ws.Range("R10").Value = 488.3132258354
$ws.Range("S10").Value = 0.07474057115881667
$ws.Range("T10").Value = 0.07474057115881669
